$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to keep text formatting so numeric-looking
# strings like "26.455.46" are not auto-converted to numbers by Excel.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "26.455.46"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.617.22"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "212.65"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "1.842.58"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "1.617.90"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "63.76"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "234.75"
$ws.Range("D18").Value = "26.459.37"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "7.73"
$ws.Range("E19").Value = "  +5.15%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0726"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "2.18"
$ws.Range("E23").Value = "  +4.34%  "
$ws.Range("D24").Value = "9.07"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "146.93"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "7.01"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "15.51"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "1.516.38"
$ws.Range("E32").Value = "  +6.85%  "
$ws.Range("D33").Value = "3.25"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").Value = "1.52"
$ws.Range("E35").Value = "  +4.09%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "0.567"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "0.830"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "5.88"
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("D43").Value = "1.755.35"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "61.59"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "0.905"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").Value = "89.77"
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("D48").Value = "1.50"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "0.0961"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "7.47"
$ws.Range("E51").Value = "  +1.19%  "

# Restore the default style on the price column so no stray
# cell-format index is left behind (matches original workbook).
$priceCol.Style = "Normal"
